$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (row 1, column A) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 31 de Marzo de 2020 a las 23:50"

# --- Update the Covid stats for the countries that received new numbers ---
# (Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#           F=Casos criticos, G=Muertes hoy, H=Muertes)

# Estados Unidos (currently row 4)
$ws.Cells.Item(4, 2).Value = 185270
$ws.Cells.Item(4, 3).Value = 21482
$ws.Cells.Item(4, 4).Value = 6347
$ws.Cells.Item(4, 5).Value = 175143
$ws.Cells.Item(4, 6).Value = 3981
$ws.Cells.Item(4, 7).Value = 639
$ws.Cells.Item(4, 8).Value = 3780

# Canada (currently row 18)
$ws.Cells.Item(18, 2).Value = 8505
$ws.Cells.Item(18, 3).Value = 1057
$ws.Cells.Item(18, 4).Value = 1162
$ws.Cells.Item(18, 5).Value = 7242
$ws.Cells.Item(18, 6).Value = 120
$ws.Cells.Item(18, 7).Value = 12
$ws.Cells.Item(18, 8).Value = 101

# Noruega (currently row 22)
$ws.Cells.Item(22, 2).Value = 4641
$ws.Cells.Item(22, 3).Value = 196
$ws.Cells.Item(22, 4).Value = 13
$ws.Cells.Item(22, 5).Value = 4589
$ws.Cells.Item(22, 6).Value = 97
$ws.Cells.Item(22, 7).Value = 7
$ws.Cells.Item(22, 8).Value = 39

# Colombia (currently row 52)
$ws.Cells.Item(52, 2).Value = 906
$ws.Cells.Item(52, 3).Value = 108
$ws.Cells.Item(52, 4).Value = 31
$ws.Cells.Item(52, 5).Value = 859
$ws.Cells.Item(52, 6).Value = 35
$ws.Cells.Item(52, 7).Value = 2
$ws.Cells.Item(52, 8).Value = 16

# Marruecos (currently row 66)
$ws.Cells.Item(66, 2).Value = 617
$ws.Cells.Item(66, 3).Value = 61
$ws.Cells.Item(66, 4).Value = 24
$ws.Cells.Item(66, 5).Value = 557
$ws.Cells.Item(66, 6).Value = 1
$ws.Cells.Item(66, 7).Value = 3
$ws.Cells.Item(66, 8).Value = 36

# Libano (currently row 71)
$ws.Cells.Item(71, 2).Value = 470
$ws.Cells.Item(71, 3).Value = 24
$ws.Cells.Item(71, 4).Value = 37
$ws.Cells.Item(71, 5).Value = 421
$ws.Cells.Item(71, 6).Value = 7
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 12

# Tunez (currently row 77) - updated numbers push it up in rank
$ws.Cells.Item(77, 2).Value = 394
$ws.Cells.Item(77, 3).Value = 32
$ws.Cells.Item(77, 4).Value = 3
$ws.Cells.Item(77, 5).Value = 381
$ws.Cells.Item(77, 6).Value = 10
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 10

# San Marino (currently row 91)
$ws.Cells.Item(91, 2).Value = 236
$ws.Cells.Item(91, 3).Value = 6
$ws.Cells.Item(91, 4).Value = 13
$ws.Cells.Item(91, 5).Value = 197
$ws.Cells.Item(91, 6).Value = 16
$ws.Cells.Item(91, 7).Value = 1
$ws.Cells.Item(91, 8).Value = 26

# Bermudas (currently row 146) - updated numbers push it up in rank
$ws.Cells.Item(146, 2).Value = 32
$ws.Cells.Item(146, 3).Value = 5
$ws.Cells.Item(146, 4).Value = 10
$ws.Cells.Item(146, 5).Value = 22
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 0

# --- Re-sort the country table (rows 4-208) descending by "Casos totales" (col B) ---
# so countries whose totals increased (Tunez, Bermudas) move up to their new rank,
# matching the way this data sheet is always kept in descending order.
$dataRange = $ws.Range("A4:H208")
$keyRange = $ws.Range("B4:B208")
$dataRange.Sort($keyRange, 2)
